$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Resultados"

# Column A width
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668

$lastRow = 104
$all = $ws.Range("A1:B$lastRow")

# Borders for the whole data range first
$all.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$all.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Header row formatting (bold white font on blue fill, centered)
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Font.Color = 16777215
$header.Interior.Color = 12419407
$header.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$header.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# Body formatting - column B (points) centered
$colB = $ws.Range("B2:B$lastRow")
$colB.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$colB.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# Turn the range into a table
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $all, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Resultados"
$tbl.TableStyle = "TableStyleMedium9"
$tbl.ShowTableStyleColumnStripes = $true
